$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 update
$ws.Range("AI3").Value = 1250

# Row 4 updates
$ws.Range("Q4").Value = 1.9
$ws.Range("R4").Value = 1.95
$ws.Range("S4").Value = 3.25
$ws.Range("T4").Value = 1.33
$ws.Range("U4").Value = 1.4
$ws.Range("V4").Value = 2.75
$ws.Range("W4").Value = 1.67
$ws.Range("X4").Value = 2.1
$ws.Range("Y4").Value = 9.5
$ws.Range("AB4").Value = 23
$ws.Range("AC4").Value = 21
$ws.Range("AJ4").Value = 10
